$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dimension-affecting data: row 2 (existing company, now id=2) ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("D2").Value = 0.0482
$ws.Range("E2").Value = -0.174
$ws.Range("F2").Value = 0.4379999999999999
$ws.Range("G2").Value = 0.3150157771468592
$ws.Range("H2").Value = 0.3150157771468592
$ws.Range("I2").Value = 0.332939218447589
$ws.Range("J2").Value = 0.2726723198355351
$ws.Range("K2").Value = 332.097
$ws.Range("L2").Value = 0.2500975246070775
$ws.Range("M2").Value = 89.148
$ws.Range("N2").Value = 0.004017231900646196
$ws.Range("O2").Value = 0.2684396426345315
$ws.Range("P2").Value = 89.148
$ws.Range("Q2").Value = 0.004017231900646196
$ws.Range("R2").Value = 0.2684396426345315
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 117.77
$ws.Range("V2").Value = 0.005307010824012907
$ws.Range("W2").Value = 0.2539386482565079
$ws.Range("X2").Value = 0.04425857756532717
$ws.Range("Y2").Value = 0.2096800706911807
$ws.Range("Z2").Value = 0.2279855640791434
$ws.Range("AA2").Value = 0.03089534142904126
$ws.Range("AB2").Value = 0.04278106991519641
$ws.Range("AC2").Value = -0.01188572848615515
$ws.Range("AD2").Value = 6624.4
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 6624.4
$ws.Range("AG2").Value = 6506.629999999999
$ws.Range("AH2").Value = 0.229887769904011
$ws.Range("AI2").Value = 0.8088696777659743
$ws.Range("AJ2").Value = 0.2267274095120815
$ws.Range("AK2").Value = 0.8060810735474663
$ws.Range("AL2").Value = 12.1
$ws.Range("AM2").Value = 12.1
$ws.Range("AN2").Value = 14.57513751375137
$ws.Range("AO2").Value = 36.53719008264463
$ws.Range("AP2").Value = 14.31601760176017
$ws.Range("AQ2").Value = 36.53719008264463

# --- Update row 3 (XP Inc.) ---
$ws.Range("A3").Value = "Brazil"
$ws.Range("B3").Value = "XP Inc. (NasdaqGS:XP)"
$ws.Range("C3").Value = "Brokerage & Investment Banking"
$ws.Range("F3").Value = 0.4379999999999999
$ws.Range("G3").Value = 0.3155314173644113
$ws.Range("H3").Value = 0.3155314173644113
$ws.Range("I3").Value = 0.3334841970279852
$ws.Range("J3").Value = 0.2703464725906693
$ws.Range("K3").Value = 331.8
$ws.Range("L3").Value = 0.2502828694274723
$ws.Range("M3").Value = 89
$ws.Range("N3").Value = 0.004013926973589262
$ws.Range("O3").Value = 0.2682338758288125
$ws.Range("P3").Value = 89
$ws.Range("Q3").Value = 0.004013926973589262
$ws.Range("R3").Value = 0.2682338758288125
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 114.4
$ws.Range("V3").Value = 0.005159474671669794
$ws.Range("W3").Value = 0.4976004799040192
$ws.Range("X3").Value = 0.04739588750110897
$ws.Range("Y3").Value = 0.4502045924029103
$ws.Range("Z3").Value = 0.2285610841005482
$ws.Range("AA3").Value = 0.06179068285808252
$ws.Range("AB3").Value = 0.04444087220084745
$ws.Range("AC3").Value = 0.01734981065723507
$ws.Range("AD3").Value = 6624.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 6624.4
$ws.Range("AG3").Value = 6510
$ws.Range("AH3").Value = 0.2300362535246482
$ws.Range("AI3").Value = 0.8110483979578094
$ws.Range("AJ3").Value = 0.2269652893022996
$ws.Range("AK3").Value = 0.8083642730309314
$ws.Range("AL3").Value = 12.1
$ws.Range("AM3").Value = 12.1
$ws.Range("AN3").Value = 14.57513751375137
$ws.Range("AO3").Value = 36.53719008264463
$ws.Range("AP3").Value = 14.32343234323432
$ws.Range("AQ3").Value = 36.53719008264463

# --- Add new row 4 (Banco Mercantil de Investimentos) ---
$ws.Range("A4").Value = "Brazil"
$ws.Range("B4").Value = "Banco Mercantil de Investimentos S.A. (BOVESPA:BMIN4)"
$ws.Range("C4").Value = "Brokerage & Investment Banking"
$ws.Range("D4").Value = 0.0482
$ws.Range("E4").Value = -0.174
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.297
$ws.Range("L4").Value = 0.1368663594470046
$ws.Range("M4").Value = 0.148
$ws.Range("N4").Value = 0.007956989247311827
$ws.Range("O4").Value = 0.4983164983164983
$ws.Range("P4").Value = 0.148
$ws.Range("Q4").Value = 0.007956989247311827
$ws.Range("R4").Value = 0.4983164983164983
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 3.37
$ws.Range("V4").Value = 0.1811827956989247
$ws.Range("W4").Value = 0.01027681660899654
$ws.Range("X4").Value = 0.04112126762954536
$ws.Range("Y4").Value = -0.03084445102054882
$ws.Range("Z4").Value = 0.089817880794702
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04112126762954536
$ws.Range("AC4").Value = -0.04112126762954536
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -3.37
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.2212738017071569
$ws.Range("AK4").Value = -0.1808910359634998
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
